$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph that follows the title.
#    (It is the 2nd paragraph in the document.)
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
if ($metaPara.Range.Text -like "Meta description*") {
    $metaPara.Range.Delete()
}

# ---------------------------------------------------------------------------
# 2) The final paragraph currently holds the (italic) image-prompt text:
#       "Create a feature image for Fire Joker ..."
#    It needs to become two paragraphs:
#       a) a new, bold paragraph containing the page title text
#          "Play Fire Joker Free - Classic Slot Game Review"
#       b) the existing paragraph (still italic) with its text replaced by
#          the meta-description sentence that used to live near the top.
# ---------------------------------------------------------------------------
$headingText = "Play Fire Joker Free - Classic Slot Game Review"
$metaDescriptionText = "Read our Fire Joker online slot game review and play it free. Learn about its features and multi-tiered bonus rounds."

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)

# Insert the heading text at the very start of the last paragraph ...
$insertPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$insertPoint.InsertBefore($headingText)

# ... then split it off into its own paragraph with a paragraph break.
$splitPos = $lastPara.Range.Start + $headingText.Length
$splitRange = $d.Range($splitPos, $splitPos)
$splitRange.InsertParagraphAfter()

# Make the newly created heading paragraph bold.
$newHeadingPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$newHeadingRange = $d.Range($newHeadingPara.Range.Start, $newHeadingPara.Range.End - 1)
$newHeadingRange.Bold = 1

# Replace the text of the remaining (italic) paragraph with the
# meta-description sentence.
$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$finalTextRange = $d.Range($finalPara.Range.Start, $finalPara.Range.End - 1)
$finalTextRange.Text = $metaDescriptionText
